$d = $word.ActiveDocument

# Update the date heading
$d.Paragraphs.Item(1).Range.Text = "2024-03-13 Wednesday"

# Update table cell values (20 rows x 5 columns)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "58+27="
$t.Cell(1,2).Range.Text = "35-27="
$t.Cell(1,3).Range.Text = "66+18="
$t.Cell(1,4).Range.Text = "7+14="
$t.Cell(1,5).Range.Text = "7+24="

$t.Cell(2,1).Range.Text = "74-39="
$t.Cell(2,2).Range.Text = "7+79="
$t.Cell(2,3).Range.Text = "58+4="
$t.Cell(2,4).Range.Text = "8+6="
$t.Cell(2,5).Range.Text = "45-8="

$t.Cell(3,1).Range.Text = "91-19="
$t.Cell(3,2).Range.Text = "24-15="
$t.Cell(3,3).Range.Text = "73-57="
$t.Cell(3,4).Range.Text = "81-28="
$t.Cell(3,5).Range.Text = "91-2="

$t.Cell(4,1).Range.Text = "3+89="
$t.Cell(4,2).Range.Text = "5+87="
$t.Cell(4,3).Range.Text = "94-87="
$t.Cell(4,4).Range.Text = "36+28="
$t.Cell(4,5).Range.Text = "62-17="

$t.Cell(5,1).Range.Text = "48+16="
$t.Cell(5,2).Range.Text = "67-49="
$t.Cell(5,3).Range.Text = "19+45="
$t.Cell(5,4).Range.Text = "69+22="
$t.Cell(5,5).Range.Text = "8+79="

$t.Cell(6,1).Range.Text = "93-47="
$t.Cell(6,2).Range.Text = "7+34="
$t.Cell(6,3).Range.Text = "23-19="
$t.Cell(6,4).Range.Text = "69+15="
$t.Cell(6,5).Range.Text = "70-2="

$t.Cell(7,1).Range.Text = "38+14="
$t.Cell(7,2).Range.Text = "28+55="
$t.Cell(7,3).Range.Text = "17+56="
$t.Cell(7,4).Range.Text = "4+28="
$t.Cell(7,5).Range.Text = "85-48="

$t.Cell(8,1).Range.Text = "81-66="
$t.Cell(8,2).Range.Text = "8+65="
$t.Cell(8,3).Range.Text = "47+16="
$t.Cell(8,4).Range.Text = "44-38="
$t.Cell(8,5).Range.Text = "16+16="

$t.Cell(9,1).Range.Text = "26+46="
$t.Cell(9,2).Range.Text = "96-27="
$t.Cell(9,3).Range.Text = "60-53="
$t.Cell(9,4).Range.Text = "90-26="
$t.Cell(9,5).Range.Text = "63+9="

$t.Cell(10,1).Range.Text = "19+49="
$t.Cell(10,2).Range.Text = "51-19="
$t.Cell(10,3).Range.Text = "47+9="
$t.Cell(10,4).Range.Text = "19+54="
$t.Cell(10,5).Range.Text = "18+39="

$t.Cell(11,1).Range.Text = "33-28="
$t.Cell(11,2).Range.Text = "82-58="
$t.Cell(11,3).Range.Text = "40-23="
$t.Cell(11,4).Range.Text = "80-54="
$t.Cell(11,5).Range.Text = "88-49="

$t.Cell(12,1).Range.Text = "72-38="
$t.Cell(12,2).Range.Text = "60-18="
$t.Cell(12,3).Range.Text = "49+36="
$t.Cell(12,4).Range.Text = "31-3="
$t.Cell(12,5).Range.Text = "27+5="

$t.Cell(13,1).Range.Text = "33-26="
$t.Cell(13,2).Range.Text = "22+9="
$t.Cell(13,3).Range.Text = "9+47="
$t.Cell(13,4).Range.Text = "55+29="
$t.Cell(13,5).Range.Text = "53-35="

$t.Cell(14,1).Range.Text = "91-84="
$t.Cell(14,2).Range.Text = "17+8="
$t.Cell(14,3).Range.Text = "55+18="
$t.Cell(14,4).Range.Text = "64-56="
$t.Cell(14,5).Range.Text = "25+56="

$t.Cell(15,1).Range.Text = "64+28="
$t.Cell(15,2).Range.Text = "79+2="
$t.Cell(15,3).Range.Text = "53+18="
$t.Cell(15,4).Range.Text = "8+4="
$t.Cell(15,5).Range.Text = "29+54="

$t.Cell(16,1).Range.Text = "39+47="
$t.Cell(16,2).Range.Text = "36+57="
$t.Cell(16,3).Range.Text = "81-44="
$t.Cell(16,4).Range.Text = "12+29="
$t.Cell(16,5).Range.Text = "81-73="

$t.Cell(17,1).Range.Text = "92-16="
$t.Cell(17,2).Range.Text = "65-37="
$t.Cell(17,3).Range.Text = "53-5="
$t.Cell(17,4).Range.Text = "21-6="
$t.Cell(17,5).Range.Text = "81-55="

$t.Cell(18,1).Range.Text = "18+75="
$t.Cell(18,2).Range.Text = "39+23="
$t.Cell(18,3).Range.Text = "68-49="
$t.Cell(18,4).Range.Text = "58+34="
$t.Cell(18,5).Range.Text = "94-26="

$t.Cell(19,1).Range.Text = "18+75="
$t.Cell(19,2).Range.Text = "64+7="
$t.Cell(19,3).Range.Text = "74-39="
$t.Cell(19,4).Range.Text = "6+77="
$t.Cell(19,5).Range.Text = "54+28="

$t.Cell(20,1).Range.Text = "31-9="
$t.Cell(20,2).Range.Text = "12+29="
$t.Cell(20,3).Range.Text = "36+55="
$t.Cell(20,4).Range.Text = "44+19="
$t.Cell(20,5).Range.Text = "64-5="
